# "Qabc i pipelines færdig" - Qa/Qb/Qc under the "L04 - Pipelines" section are
# now done: clear the leftover "mangler ..." remark text from Qb (C12) and
# Qc (C13) and mark them with the same green "done" fill used for the other
# finished cells (C7, C8, C11, ...). Also move the active selection to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doneGreen = 5287936  # RGB(0, 176, 80) == fill color used for completed items

$ws.Range("C12").ClearContents()
$ws.Range("C12").Interior.Color = $doneGreen

$ws.Range("C13").ClearContents()
$ws.Range("C13").Interior.Color = $doneGreen

$ws.Range("C8").Select() | Out-Null
